$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2, pushing all existing data
# rows (2..37) down to (3..38).
$ws.Range("A2").EntireRow.Insert()

# The inserted row copied formatting from the header row (bold, borders).
# Clear that so it matches the plain formatting used by the other data rows.
$ws.Range("A2:R2").ClearFormats()

# Populate the new row 2 with the latest weekly price record.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").Value = 44699
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 100112052
$ws.Range("G2").Value = "Albahaca"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = 2250
$ws.Range("N2").Value = "$/paquete"
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 2250
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = "Hortaliza"
